# Apply the "Deploying to gh-pages" metadata refresh:
#  - Alvearie/IBM -> LinuxForHealth rebrand of the canonical URL + publisher
#  - bump version 7.0.0 -> 8.0.0
#  - bump the generation timestamp
#  - clear the now-redundant Constraint(s) text that used to be duplicated
#    on the root "Extension" row (it now only lives on "Element.extension")

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/patient-citizenship"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
# The extension's canonical URL is mirrored here as the fixed value of Extension.url
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/patient-citizenship"
# Constraint text is now only shown on the Element.extension row, not on Extension itself
$elements.Range("AI2").Value = ""
